# Historia_Grafica.docx - "Feito o corpo da primeira pagina" edit
#
# 1. "plasticos, etc...) Atraves ... off-set, digital, rotogravura,"
#    -> "plasticos etc.) Atraves ... offset, digital, fotogravura,"
#    (fixes punctuation/typos, also clears the spelling/grammar proofErr
#    markers Word had recorded over that stretch of text)
# 2. Hyperlink over "cartoes de visita" no longer carries an (empty)
#    ScreenTip / tooltip.
# 3. Minor proofing-marker cleanup around "flyers" and "Johannes Gutenberg"
#    (no visible text change there).

$d = $word.ActiveDocument

# --- 1. substrate list + printing techniques -------------------------------
$d.Content.Find.Execute(
    "plásticos, etc...) Através de um sistema de impressão, como off-set, digital, rotogravura,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "plásticos etc.) Através de um sistema de impressão, como offset, digital, fotogravura,",
    2) | Out-Null

# --- 2. tidy the "flyers" mention (formatting/proofing only) ---------------
$d.Content.Find.Execute(
    "flyers, cartazes, entre outros",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "flyers, cartazes, entre outros",
    2) | Out-Null

# --- 3. tidy "Johannes Gutenberg" (formatting/proofing only) ---------------
$d.Content.Find.Execute(
    "Johannes Gutenberg",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Johannes Gutenberg",
    2) | Out-Null

# --- 4. drop the (empty) tooltip/ScreenTip on the "cartões de visita" link -
foreach ($hl in $d.Hyperlinks) {
    $hl.ScreenTip = ""
}
